$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''61.660.84'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '''2.896.17'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.17%  '
$ws.Range('D4').Value = '''0.997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').Value = '''584.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').Value = '''145.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.08%  '
$ws.Range('D8').Value = '''0.500'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.30%  '
$ws.Range('D9').Value = '''2.889.44'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').Value = '''7.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.18%  '
$ws.Range('D11').Value = '''0.148'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +7.23%  '
$ws.Range('D12').Value = '''0.432'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').Value = '''0.0000233'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +7.00%  '
$ws.Range('D14').Value = '''31.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('E15').Value = '  -1.07%  '
$ws.Range('D16').Value = '''3.379.73'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').Value = '''61.639.38'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '''6.55'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.24%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '''2.895.30'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.48%  '
$ws.Range('D20').Value = '''430.18'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.34%  '
$ws.Range('D21').Value = '''13.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.21%  '
$ws.Range('D22').Value = '''0.654'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('D23').Value = '''6.83'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.38%  '
$ws.Range('D24').Value = '''79.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Value = '''10.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.05%  '
$ws.Range('D26').Value = '''11.80'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.05%  '
$ws.Range('D27').Value = '''2.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('D29').Value = '''7.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.26%  '
$ws.Range('D30').Value = '''2.55'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = '''0.0000101'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +19.79%  '
$ws.Range('D32').Value = '''2.10'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.02%  '
$ws.Range('D33').Value = '''0.107'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('D34').Value = '''0.997'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('D35').Value = '''25.73'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '''0.973'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.35%  '
$ws.Range('D37').Value = '''3.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +9.10%  '
$ws.Range('D38').Value = '''5.46'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '''49.10'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.32%  '
$ws.Range('D40').Value = '''1.97'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.46%  '
$ws.Range('D41').Value = '''8.29'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.85%  '
$ws.Range('D42').Value = '''0.114'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.50%  '
$ws.Range('D43').Value = '''0.271'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('D44').Value = '''39.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.20%  '
$ws.Range('D45').Value = '''134.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('D46').Value = '''2.671.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.08%  '
$ws.Range('D47').Value = '''0.0335'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.19%  '
$ws.Range('D48').Value = '''348.52'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.03%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  +0.86%  '
$ws.Range('D51').Value = '''22.29'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.29%  '
